$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "c20cd1810746ebf932bb57d9f8f33030"
$ws.Range("B44").Value = "64d97a2435ca528474a9ee1b62a5969d"
$ws.Range("B54").Value = "352a93de4c96c92f8d698df76762f5fa"
$ws.Range("B74").Value = "7ab7fef2fd4db72bbdb0720aafcbc718"
$ws.Range("B89").Value = "677808ed7f974be62cdfb69b4daed013"
$ws.Range("B99").Value = "7295799e6758bfbfe9f01c6adf1aca08"
$ws.Range("B110").Value = "1bd5e3b761a52acf1e20b0c69324b0d1"
$ws.Range("B126").Value = "30992a194a56e3775d7bc9fa5a64bc24"
$ws.Range("B161").Value = "10f1715cd7ab53d5a3f38c26ac1e512f"
$ws.Range("B168").Value = "b59d55c420b531bf2814747715b21456"
$ws.Range("B191").Value = "acaccb83a22399e165da21f829a51351"
$ws.Range("B198").Value = "386b9f534bea2b82a41071574dc0218f"
$ws.Range("B222").Value = "71b88b81eea7e7cd24ef136b4f66f21b"
$ws.Range("B227").Value = "79d7ac27c02b8ee4b146a8ebaf9cdac1"
$ws.Range("B229").Value = "03ddbb620ca14ff25340edfbe05fe1de"
$ws.Range("B232").Value = "ae22bcdb5a3d16e8e1bb7667b80435a8"
$ws.Range("B278").Value = "6ca2b727497da9da297e10d0e74f11fc"
$ws.Range("B345").Value = "1d0565d3900a06151050ed3f0730ef7c"
$ws.Range("B419").Value = "2f36e7fae61a39e97cd825cd8a551d49"
$ws.Range("B478").Value = "19b25a4ce25f6f97839a85d363ab88b0"
$ws.Range("B480").Value = "76e3d164f5a0404b0df223484a58660f"
$ws.Range("B501").Value = "2f3dfc70d7f041da9765e62f76ca913a"
$ws.Range("B726").Value = "63c9f9c955a1cd66bf998e68d6445a72"
$ws.Range("B733").Value = "4c378edcdadf5352ae31165b2ead8eaa"
$ws.Range("B768").Value = "fa3438559eb36bcd278952eeb9ffd616"
$ws.Range("B816").Value = "e156ff61a68c1b859d559b0ba2bd01c0"
$ws.Range("B825").Value = "74f20965bca711405d4b5008fd88b85c"
$ws.Range("B827").Value = "7c0d8b2c888ea89da57dac14fe891e28"
$ws.Range("B850").Value = "ee5f9b6f034b61262ef8922f4d4f5ebd"
$ws.Range("B862").Value = "15adcc8626573003a2667afe259f8d2e"
